$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price strings so they remain text
# (matches source data, which stores prices as text, e.g. "231.58")
foreach ($addr in @("D5","D7","D10","D11","D12","D14","D15","D16","D17","D21","D22","D23","D24","D27","D28","D30","D31","D32","D34","D36","D37","D38","D43","D45","D46","D49","D47")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '41.729.67'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '2.232.56'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '231.58'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("E6").Value = '  -2.23%  '
$ws.Range("D7").Value = '60.21'
$ws.Range("E7").Value = '  -7.01%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("E9").Value = '  -2.11%  '
$ws.Range("D10").Value = '57.97'
$ws.Range("E10").Value = '  -2.05%  '
$ws.Range("D11").Value = '0.0898'
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").Value = '0.104'
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("D13").Value = '2.565.11'
$ws.Range("E13").Value = '  -0.89%  '
$ws.Range("D14").Value = '15.47'
$ws.Range("E14").Value = '  -4.91%  '
$ws.Range("D15").Value = '22.52'
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = '5.63'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '0.800'
$ws.Range("E17").Value = '  -4.45%  '
$ws.Range("D18").Value = '2.234.29'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").Value = '41.682.28'
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("D20").Value = '0.0₃0910'
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("D21").Value = '72.52'
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("D22").Value = '6.14'
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("D23").Value = '247.91'
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("D24").Value = '0.998'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  -1.42%  '
$ws.Range("E26").Value = '  -1.04%  '
$ws.Range("D27").Value = '9.77'
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("D28").Value = '169.43'
$ws.Range("E28").Value = '  -2.23%  '
$ws.Range("E29").Value = '  -2.59%  '
$ws.Range("D30").Value = '19.90'
$ws.Range("E30").Value = '  -2.81%  '
$ws.Range("D31").Value = '1.42'
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("D32").Value = '2.56'
$ws.Range("E32").Value = '  -8.87%  '
$ws.Range("E33").Value = '  -2.10%  '
$ws.Range("D34").Value = '5.02'
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("E35").Value = '  -1.50%  '
$ws.Range("D36").Value = '0.0655'
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("D37").Value = '6.55'
$ws.Range("E37").Value = '  -8.87%  '
$ws.Range("D38").Value = '2.41'
$ws.Range("E38").Value = '  -2.13%  '
$ws.Range("E39").Value = '  -8.21%  '
$ws.Range("E40").Value = '  +1.60%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("E42").Value = '  +0.20%  '
$ws.Range("D43").Value = '8.68'
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("D45").Value = '99.00'
$ws.Range("E45").Value = '  -3.49%  '
$ws.Range("D46").Value = '0.0963'
$ws.Range("E46").Value = '  +1.89%  '
$ws.Range("D49").Value = '16.62'
$ws.Range("E49").Value = '  -7.27%  '
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("E51").Value = '  -3.05%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '4.40'
$ws.Range("E47").Value = '  -9.47%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.477.20'
$ws.Range("E48").Value = '  -2.54%  '
